$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading percent values for rows 2-25 (case with 380 kV done)
$data = @{
    2 = @{ "B"="21.41559893994325"; "C"="7.389307731227945"; "D"="9.424764307976707"; "F"="50.31013185936422"; "G"="3.705335002253343"; "L"="10.85251586285393" }
    3 = @{ "B"="21.10440418298965"; "C"="6.962254345895776"; "D"="9.300060236560098"; "F"="49.11768912723406"; "G"="3.710472821973509"; "L"="10.83561971778492" }
    4 = @{ "B"="20.92104351675718"; "C"="6.68521250824555"; "D"="9.222972015959229"; "F"="48.37893750521798"; "G"="3.71378327579222"; "L"="10.82771821263354" }
    5 = @{ "B"="20.84837192724169"; "C"="6.568603389908654"; "D"="9.191444315663141"; "F"="48.07659277075595"; "G"="3.71517167943758"; "L"="10.82512078373352" }
    6 = @{ "B"="20.83643174067576"; "C"="6.549016929399793"; "D"="9.186202796607809"; "F"="48.02632102676768"; "G"="3.715404605928421"; "L"="10.8247270907922" }
    7 = @{ "B"="20.92005500347865"; "C"="6.683654881513862"; "D"="9.222547260079818"; "F"="48.37486474829282"; "G"="3.713801840645847"; "L"="10.82768066175083" }
    8 = @{ "B"="21.30676517997379"; "C"="7.245145666123154"; "D"="9.381884355235371"; "F"="49.90055593650575"; "G"="3.707074296563738"; "L"="10.84617739915372" }
    9 = @{ "B"="22.12127160582199"; "C"="8.228125884841335"; "D"="9.689479235354325"; "F"="52.82513092391925"; "G"="3.695109361277635"; "L"="10.90202017445623" }
    10 = @{ "B"="22.74697585181093"; "C"="8.878071681434282"; "D"="9.911519848063703"; "F"="54.91384621678791"; "G"="3.687055212831933"; "L"="10.9548841003453" }
    11 = @{ "B"="23.03605157945968"; "C"="9.158078609722054"; "D"="10.01146958991017"; "F"="55.84740089687759"; "G"="3.683548511352695"; "L"="10.98147323959016" }
    12 = @{ "B"="23.14603167399187"; "C"="9.261862624655276"; "D"="10.04914896439901"; "F"="56.19824861525474"; "G"="3.682243013847683"; "L"="10.99190386278359" }
    13 = @{ "B"="23.122324623496"; "C"="9.239610817040681"; "D"="10.04104184337137"; "F"="56.12281023457982"; "G"="3.682523182191919"; "L"="10.98964140703714" }
    14 = @{ "B"="23.04509012622558"; "C"="9.166662023342784"; "D"="10.01457298791613"; "F"="55.87632037529333"; "G"="3.683440659144535"; "L"="10.98232414213211" }
    15 = @{ "B"="22.99784495616813"; "C"="9.121686075802433"; "D"="9.998337478865924"; "F"="55.72498265140823"; "G"="3.684005553712108"; "L"="10.9778891230765" }
    16 = @{ "B"="22.72816267036053"; "C"="8.859457716425959"; "D"="9.904965123059119"; "F"="54.85247749707298"; "G"="3.687287530749646"; "L"="10.95319728120935" }
    17 = @{ "B"="22.56376878232199"; "C"="8.694580436308019"; "D"="9.847400849677692"; "F"="54.31276291581884"; "G"="3.689341039897002"; "L"="10.93869811920634" }
    18 = @{ "B"="22.46963902709971"; "C"="8.598273385326054"; "D"="9.814192863608769"; "F"="54.00078517048059"; "G"="3.690536968927503"; "L"="10.9305980033192" }
    19 = @{ "B"="22.43784519505881"; "C"="8.565412258554923"; "D"="9.802932821334535"; "F"="53.89489784758358"; "G"="3.690944438380001"; "L"="10.92789666281293" }
    20 = @{ "B"="22.58122569273089"; "C"="8.712284443310413"; "D"="9.853538990166308"; "F"="54.37037878372639"; "G"="3.689120909444642"; "L"="10.94021682476913" }
    21 = @{ "B"="23.06776282105478"; "C"="9.188149843924625"; "D"="10.02235226642723"; "F"="55.9487949928249"; "G"="3.683170567098816"; "L"="10.98446360883237" }
    22 = @{ "B"="23.38867937025589"; "C"="9.48605191488636"; "D"="10.13168626170452"; "F"="56.96470799662045"; "G"="3.679412234673094"; "L"="11.01548887469154" }
    23 = @{ "B"="23.21717210116512"; "C"="9.328253485339287"; "D"="10.07342932312227"; "F"="56.42401740034112"; "G"="3.681406242989596"; "L"="10.99873855608915" }
    24 = @{ "B"="22.57333222267904"; "C"="8.704285186765611"; "D"="9.850764289115915"; "F"="54.3443358880933"; "G"="3.689220382629695"; "L"="10.93952948374048" }
    25 = @{ "B"="21.89562463061479"; "C"="7.97488774507755"; "D"="9.606889591184755"; "F"="52.04318700400413"; "G"="3.698215991006295"; "L"="10.88482584843357" }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = [double]$rowData[$col]
    }
}
